$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spot_PT")

$ws.Range("A2").Value = 46014
$ws.Range("B2").Value = 70.89
$ws.Range("C2").Value = 62.05
$ws.Range("D2").Value = 61.35
$ws.Range("E2").Value = 59.43
$ws.Range("F2").Value = 56.65
$ws.Range("G2").Value = 60.85
$ws.Range("H2").Value = 70.40000000000001
$ws.Range("I2").Value = 81
$ws.Range("J2").Value = 87.17
$ws.Range("K2").Value = 88.77
$ws.Range("L2").Value = 85.14
$ws.Range("M2").Value = 83.40000000000001
$ws.Range("N2").Value = 84.18000000000001
$ws.Range("O2").Value = 81.17
$ws.Range("P2").Value = 80.67
$ws.Range("Q2").Value = 84.31999999999999
$ws.Range("R2").Value = 82.87
$ws.Range("S2").Value = 85.89
$ws.Range("T2").Value = 87.45999999999999
$ws.Range("U2").Value = 86.09
$ws.Range("V2").Value = 89.97
$ws.Range("W2").Value = 89.06999999999999
$ws.Range("X2").Value = 81.37
$ws.Range("Y2").Value = 69.81999999999999
$ws.Range("Z2").Value = 77.92
$ws.Range("AA2").Value = "8h-12h"
$ws.Range("AB2").Value = 86.12
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 89.52
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 87.97
$ws.Range("AG2").Value = "0h-23h"
